$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: T5YIFR ---
$ws.Range("N29").Value = 46076
$ws.Range("Q29").Value = 2.12
$ws.Range("R29").Value = 2.13
$ws.Range("S29").Value = 2.15
$ws.Range("T29").Value = 2.15
$ws.Range("U29").Value = 2.13

# --- Row 30: T10YIE ---
$ws.Range("N30").Value = 46076
$ws.Range("Q30").Value = 2.26
$ws.Range("R30").Value = 2.28
$ws.Range("S30").Value = 2.29
$ws.Range("T30").Value = 2.29
$ws.Range("U30").Value = 2.26

# --- Row 37: CSUSHPINSA (M/M) ---
# Date cell also needs the highlighted-date style used by N29 (style index 48).
$ws.Range("N29").Copy()
$ws.Range("N37").PasteSpecial(-4122)
$ws.Range("N37").Value = 45992
$ws.Range("Q37").Value = -0.00273486156666769
$ws.Range("R37").Value = -0.0006239385435752309
$ws.Range("S37").Value = -0.001255429809922437
$ws.Range("T37").Value = -0.002710199627120158
$ws.Range("U37").Value = -0.003326625695690755

# --- Row 38: CSUSHPINSA (Y/Y) ---
$ws.Range("N29").Copy()
$ws.Range("N38").PasteSpecial(-4122)
$ws.Range("N38").Value = 45992
$ws.Range("Q38").Value = 0.01272340965111129
$ws.Range("R38").Value = 0.01427415269420379
$ws.Range("S38").Value = 0.01402408537902678
$ws.Range("T38").Value = 0.01316616157883059
$ws.Range("U38").Value = 0.01475691239986718

$excel.CutCopyMode = $false

# --- Row 47: FFR (date only, no value changes) ---
$ws.Range("N47").Value = 46073

# --- Row 48: 2y UST ---
$ws.Range("N48").Value = 46073
$ws.Range("Q48").Value = 3.48
$ws.Range("R48").Value = 3.47
$ws.Range("S48").Value = 3.47
$ws.Range("T48").Value = 3.43
$ws.Range("U48").Value = 3.4

# --- Row 49: 5y UST ---
$ws.Range("N49").Value = 46073
$ws.Range("Q49").Value = 3.65
$ws.Range("R49").Value = 3.65
$ws.Range("S49").Value = 3.66
$ws.Range("T49").Value = 3.63
$ws.Range("U49").Value = 3.61

# --- Row 50: 10y UST ---
$ws.Range("N50").Value = 46073
$ws.Range("Q50").Value = 4.08
$ws.Range("R50").Value = 4.08
$ws.Range("S50").Value = 4.09
$ws.Range("T50").Value = 4.05
$ws.Range("U50").Value = 4.04

# --- Row 52: BAA ---
$ws.Range("N52").Value = 46073
$ws.Range("Q52").Value = 5.77
$ws.Range("R52").Value = 5.76
$ws.Range("S52").Value = 5.76
$ws.Range("T52").Value = 5.75
$ws.Range("U52").Value = 5.76
